$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3759955672705644
$ws.Range("C2").Value = 0.05892570557490728
$ws.Range("D2").Value = 0.3019038439208259
$ws.Range("F2").Value = 1.435412107694724
$ws.Range("G2").Value = 0.00245950757720754
$ws.Range("I2").Value = 0.6733453806067367
$ws.Range("J2").Value = 0.3500415298297384
$ws.Range("K2").Value = 0.4073226755704411
$ws.Range("O2").Value = 3.24403693567946
$ws.Range("B3").Value = 0.3357964107202918
$ws.Range("C3").Value = 0.05171371804611624
$ws.Range("D3").Value = 0.2903279750729126
$ws.Range("F3").Value = 1.438838889608803
$ws.Range("G3").Value = 0.002461975659394211
$ws.Range("I3").Value = 0.6809936148411602
$ws.Range("J3").Value = 0.3387173486612767
$ws.Range("K3").Value = 0.3620647754082427
$ws.Range("O3").Value = 3.269168133293618
$ws.Range("B4").Value = 0.31111549967261
$ws.Range("C4").Value = 0.04727583534052826
$ws.Range("D4").Value = 0.2833466330215799
$ws.Range("F4").Value = 1.441765087337522
$ws.Range("G4").Value = 0.002463571694942387
$ws.Range("I4").Value = 0.6860660211097596
$ws.Range("J4").Value = 0.3319685340679683
$ws.Range("K4").Value = 0.334262452320246
$ws.Range("O4").Value = 3.286474590845188
$ws.Range("B5").Value = 0.3010588307464559
$ws.Range("C5").Value = 0.04546500888464777
$ws.Range("D5").Value = 0.2805335761227212
$ws.Range("F5").Value = 1.443164246531055
$ws.Range("G5").Value = 0.002464242425721116
$ws.Range("I5").Value = 0.6882276527512587
$ws.Range("J5").Value = 0.3292696841089509
$ws.Range("K5").Value = 0.3229299683888485
$ws.Range("O5").Value = 3.293998566925865
$ws.Range("B6").Value = 0.2993890073194336
$ws.Range("C6").Value = 0.04516418289485102
$ws.Range("D6").Value = 0.2800684007952583
$ws.Range("F6").Value = 1.443409059936535
$ws.Range("G6").Value = 0.002464355029945451
$ws.Range("I6").Value = 0.688592301160984
$ws.Range("J6").Value = 0.328824642762271
$ws.Range("K6").Value = 0.3210480676326313
$ws.Range("O6").Value = 3.295276385338951
$ws.Range("B7").Value = 0.3109798670386681
$ws.Range("C7").Value = 0.04725142331700738
$ws.Range("D7").Value = 0.2833085658175349
$ws.Range("F7").Value = 1.441783119960981
$ws.Range("G7").Value = 0.002463580658138041
$ws.Range("I7").Value = 0.6860947907317403
$ws.Range("J7").Value = 0.3319319285586033
$ws.Range("K7").Value = 0.3341096289905749
$ws.Range("O7").Value = 3.286574153230575
$ws.Range("B8").Value = 0.3621349367290065
$ws.Range("C8").Value = 0.05644107231395878
$ws.Range("D8").Value = 0.2978863532768798
$ws.Range("F8").Value = 1.436423009504907
$ws.Range("G8").Value = 0.002460341875055765
$ws.Range("I8").Value = 0.6759043464190491
$ws.Range("J8").Value = 0.3460945374036157
$ws.Range("K8").Value = 0.3917209922906011
$ws.Range("O8").Value = 3.252312709341851
$ws.Range("B9").Value = 0.4624399332225835
$ws.Range("C9").Value = 0.07438221404683532
$ws.Range("D9").Value = 0.3274709397742583
$ws.Range("F9").Value = 1.432437807328967
$ws.Range("G9").Value = 0.002454627583894626
$ws.Range("I9").Value = 0.658909520691914
$ws.Range("J9").Value = 0.3754910454439653
$ws.Range("K9").Value = 0.504564069006932
$ws.Range("O9").Value = 3.20002149152333
$ws.Range("B10").Value = 0.5361056119848513
$ws.Range("C10").Value = 0.08751246984505201
$ws.Range("D10").Value = 0.3498114528719327
$ws.Range("F10").Value = 1.433494216346944
$ws.Range("G10").Value = 0.002450813733682289
$ws.Range("I10").Value = 0.648248533455142
$ws.Range("J10").Value = 0.3980851130358332
$ws.Range("K10").Value = 0.5873660660144253
$ws.Range("O10").Value = 3.170701887096897
$ws.Range("B11").Value = 0.5696076171546736
$ws.Range("C11").Value = 0.09347422378854731
$ws.Range("D11").Value = 0.3601054679795368
$ws.Range("F11").Value = 1.434841397810573
$ws.Range("G11").Value = 0.00244916136939555
$ws.Range("I11").Value = 0.6437955720218582
$ws.Range("J11").Value = 0.4085816897083987
$ws.Range("K11").Value = 0.6250079711611249
$ws.Range("O11").Value = 3.159343418386214
$ws.Range("B12").Value = 0.5822921589015664
$ws.Range("C12").Value = 0.09573009884158523
$ws.Range("D12").Value = 0.3640223053939735
$ws.Range("F12").Value = 1.435476235604568
$ws.Range("G12").Value = 0.002448547474309909
$ws.Range("I12").Value = 0.6421664718513327
$ws.Range("J12").Value = 0.4125879461750799
$ws.Range("K12").Value = 0.6392578211799105
$ws.Range("O12").Value = 3.155327180454208
$ws.Range("B13").Value = 0.5795604127167451
$ws.Range("C13").Value = 0.09524433337426785
$ws.Range("D13").Value = 0.363177914301275
$ws.Range("F13").Value = 1.435333965205558
$ws.Range("G13").Value = 0.002448679162723745
$ws.Range("I13").Value = 0.6425147850824686
$ws.Range("J13").Value = 0.4117237289343336
$ws.Range("K13").Value = 0.6361890644939479
$ws.Range("O13").Value = 3.156179469318801
$ws.Range("B14").Value = 0.570651224718091
$ws.Range("C14").Value = 0.09365985081814188
$ws.Range("D14").Value = 0.3604273342849922
$ws.Range("F14").Value = 1.434891126931319
$ws.Range("G14").Value = 0.002449110627316997
$ws.Range("I14").Value = 0.6436603993254089
$ws.Range("J14").Value = 0.4089106567980849
$ws.Range("K14").Value = 0.6261804055699258
$ws.Range("O14").Value = 3.159007286804695
$ws.Range("B15").Value = 0.5651938174713393
$ws.Range("C15").Value = 0.09268908385993768
$ws.Range("D15").Value = 0.3587449577085806
$ws.Range("F15").Value = 1.434636116043777
$ws.Range("G15").Value = 0.002449376449313753
$ws.Range("I15").Value = 0.6443695648915977
$ws.Range("J15").Value = 0.4071916627307104
$ws.Range("K15").Value = 0.6200492291179955
$ws.Range("O15").Value = 3.160776528715758
$ws.Range("B16").Value = 0.5339159416627695
$ws.Range("C16").Value = 0.08712262065571963
$ws.Range("D16").Value = 0.3491413421390632
$ws.Range("F16").Value = 1.433423619488096
$ws.Range("G16").Value = 0.002450923376355059
$ws.Range("I16").Value = 0.6485475362591266
$ws.Range("J16").Value = 0.3974035341576752
$ws.Range("K16").Value = 0.584905512570856
$ws.Range("O16").Value = 3.17148404692432
$ws.Range("B17").Value = 0.5147252302097343
$ws.Range("C17").Value = 0.08370482496144405
$ws.Range("D17").Value = 0.3432833329967764
$ws.Range("F17").Value = 1.432901804009219
$ws.Range("G17").Value = 0.002451893476368994
$ws.Range("I17").Value = 0.6512122696955096
$ws.Range("J17").Value = 0.3914547931897516
$ws.Range("K17").Value = 0.5633390376632121
$ws.Range("O17").Value = 3.178559935122394
$ws.Range("B18").Value = 0.5036864440430122
$ws.Range("C18").Value = 0.08173794439279902
$ws.Range("D18").Value = 0.3399263195322817
$ws.Range("F18").Value = 1.432683235299635
$ws.Range("G18").Value = 0.002452459227979203
$ws.Range("I18").Value = 0.6527822952215487
$ws.Range("J18").Value = 0.3880537923747767
$ws.Range("K18").Value = 0.5509322476375473
$ws.Range("O18").Value = 3.182816065517869
$ws.Range("B19").Value = 0.499948784179054
$ws.Range("C19").Value = 0.08107181405721064
$ws.Range("D19").Value = 0.3387918195898578
$ws.Range("F19").Value = 1.432623238512704
$ws.Range("G19").Value = 0.002452652118731612
$ws.Range("I19").Value = 0.6533202897433874
$ws.Range("J19").Value = 0.3869058020545566
$ws.Range("K19").Value = 0.546731144166074
$ws.Range("O19").Value = 3.184289098465399
$ws.Range("B20").Value = 0.5167682014573813
$ws.Range("C20").Value = 0.08406876499142868
$ws.Range("D20").Value = 0.3439056501011919
$ws.Range("F20").Value = 1.432948909889035
$ws.Range("G20").Value = 0.002451789403235763
$ws.Range("I20").Value = 0.6509247389752204
$ws.Range("J20").Value = 0.3920859190225485
$ws.Range("K20").Value = 0.5656350714029372
$ws.Range("O20").Value = 3.177787414184934
$ws.Range("B21").Value = 0.573268127131854
$ws.Range("C21").Value = 0.09412529894007093
$ws.Range("D21").Value = 0.3612347395213646
$ws.Range("F21").Value = 1.435017814668441
$ws.Range("G21").Value = 0.002448983575443187
$ws.Range("I21").Value = 0.6433223532289105
$ws.Range("J21").Value = 0.4097360710833584
$ws.Range("K21").Value = 0.6291203142452559
$ws.Range("O21").Value = 3.1581689514083
$ws.Range("B22").Value = 0.6101824781230221
$ws.Range("C22").Value = 0.1006878053906632
$ws.Range("D22").Value = 0.3726693305326023
$ws.Range("F22").Value = 1.437096781423719
$ws.Range("G22").Value = 0.002447218672056063
$ws.Range("I22").Value = 0.638686847405765
$ws.Range("J22").Value = 0.4214546775760226
$ws.Range("K22").Value = 0.6705860692972578
$ws.Range("O22").Value = 3.147008330825315
$ws.Range("B23").Value = 0.5904818756000907
$ws.Range("C23").Value = 0.09718621922451121
$ws.Range("D23").Value = 0.3665565473622507
$ws.Range("F23").Value = 1.435920670518428
$ws.Range("G23").Value = 0.002448154350818158
$ws.Range("I23").Value = 0.6411303983731678
$ws.Range("J23").Value = 0.4151834657880329
$ws.Range("K23").Value = 0.6484575836000204
$ws.Range("O23").Value = 3.152812845360245
$ws.Range("B24").Value = 0.5158445921435941
$ws.Range("C24").Value = 0.08390423375379896
$ws.Range("D24").Value = 0.3436242667814042
$ws.Range("F24").Value = 1.432927359674835
$ws.Range("G24").Value = 0.00245183642965313
$ws.Range("I24").Value = 0.6510546131426658
$ws.Range("J24").Value = 0.3918005278254952
$ws.Range("K24").Value = 0.5645970592175615
$ws.Range("O24").Value = 3.178136084975591
$ws.Range("B25").Value = 0.4353081599592485
$ws.Range("C25").Value = 0.06953745554552881
$ws.Range("D25").Value = 0.3193610822422386
$ws.Range("F25").Value = 1.432816646127776
$ws.Range("G25").Value = 0.002456105661860402
$ws.Range("I25").Value = 0.663186792600321
$ws.Range("J25").Value = 0.3673639994935485
$ws.Range("K25").Value = 0.4740536121243508
$ws.Range("O25").Value = 3.212571124365951
